$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ---------------------------------------------------------------------------
# 1. Insert a new column before Z. This shifts the old Z:AE (web, webalert,
#    webcookie, ws, ws.async, xml) one column to the right, becoming AA:AF,
#    and leaves a blank column Z for the new "tn.5250" category.
# ---------------------------------------------------------------------------
$ws.Columns("Z:Z").Insert(-4161)

# Populate the new "tn.5250" column (header + 5 command rows).
$ws.Cells.Item(1, 26).Value = "tn.5250"
$ws.Cells.Item(2, 26).Value = "close(profile)"
$ws.Cells.Item(3, 26).Value = "open(profile)"
$ws.Cells.Item(4, 26).Value = "saveText(profile,var)"
$ws.Cells.Item(5, 26).Value = "typeKeys(profile,keystrokes)"
$ws.Cells.Item(6, 26).Value = "updateScreenFields(profile)"

# ---------------------------------------------------------------------------
# 2. Shift column A (the "target" index) rows 26-31 down to 27-32, then
#    insert the new "tn.5250" entry at A26 (keeps alphabetical order).
#    Done manually (read-then-write) because a single-cell Range.Insert on
#    this runtime shifts the whole row instead of just the one column.
# ---------------------------------------------------------------------------
$colA = @{}
for ($r = 26; $r -le 31; $r++) {
    $colA[$r] = $ws.Cells.Item($r, 1).Value()
}
for ($r = 26; $r -le 31; $r++) {
    $ws.Cells.Item($r + 1, 1).Value = $colA[$r]
}
$ws.Cells.Item(26, 1).Value = "tn.5250"

# ---------------------------------------------------------------------------
# 3. Image commands (column K): rename colorbit's first param, and insert
#    the new ocr(image,saveVar) command between crop(...) and resize(...).
# ---------------------------------------------------------------------------
$ws.Cells.Item(2, 11).Value = "colorbit(image,bit,saveTo)"

$colK = @{}
for ($r = 6; $r -le 7; $r++) {
    $colK[$r] = $ws.Cells.Item($r, 11).Value()
}
for ($r = 6; $r -le 7; $r++) {
    $ws.Cells.Item($r + 1, 11).Value = $colK[$r]
}
$ws.Cells.Item(6, 11).Value = "ocr(image,saveVar)"

# ---------------------------------------------------------------------------
# 4. Update the defined names that refer to ranges shifted above, and add
#    the brand new "tn.5250" named range.
# ---------------------------------------------------------------------------
$wb.Names.Item("image").RefersTo = "='#system'!`$K`$2:`$K`$8"
$wb.Names.Item("target").RefersTo = "='#system'!`$A`$2:`$A`$32"
$wb.Names.Item("web").RefersTo = "='#system'!`$AA`$2:`$AA`$144"
$wb.Names.Item("webalert").RefersTo = "='#system'!`$AB`$2:`$AB`$8"
$wb.Names.Item("webcookie").RefersTo = "='#system'!`$AC`$2:`$AC`$10"
$wb.Names.Item("ws").RefersTo = "='#system'!`$AD`$2:`$AD`$17"
$wb.Names.Item("ws.async").RefersTo = "='#system'!`$AE`$2:`$AE`$8"
$wb.Names.Item("xml").RefersTo = "='#system'!`$AF`$2:`$AF`$27"
$wb.Names.Add("tn.5250", "='#system'!`$Z`$2:`$Z`$6")
